$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "Save" header in H1, copying the existing header style (G1) so
# the new cell gets the same bold/border/center formatting as the rest of
# row 1 instead of creating a brand new style entry.
$ws.Range("G1").Copy($ws.Range("H1")) | Out-Null
$ws.Range("H1").Value = "Save"

# Fill in the "Save" indicator values for each data row (1 = saved game).
$saveValues = @{
    2  = 0
    3  = 0
    4  = 0
    5  = 0
    6  = 0
    7  = 1
    8  = 1
    9  = 0
    10 = 1
    11 = 0
    12 = 0
    13 = 1
    14 = 0
    15 = 1
    16 = 1
}

foreach ($row in $saveValues.Keys) {
    $ws.Range("H$row").Value = $saveValues[$row]
}

Write-Host "Added Save column (H1:H16)"
